$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 duplicates row 2's data (venue/date/result/teams/batsman + stats)
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 04 2020"
$ws.Range("C3").Value = "Super Kings won by 10 wickets (with 14 balls remaining)"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Nicholas Pooran †"

# These look numeric, so prefix with an apostrophe to force them to stay
# text (matching the source file's t="str" cells) instead of becoming
# real numbers.
$ws.Range("G3").Value = "'33"
$ws.Range("H3").Value = "'17"
$ws.Range("I3").Value = "'1"
$ws.Range("J3").Value = "'3"
$ws.Range("K3").Value = "'194.11"
